$p = $ppt.ActivePresentation
$layout = $p.SlideMaster.CustomLayouts.Item(1)
$s = $p.Slides.AddSlide(1, $layout)
